# Apply updated crypto price/volume data (and two name/link/price swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '28.929.10'

Set-TextCell 3 4 '1.910.02'
Set-TextCell 3 5 '  -1.82%  '

Set-TextCell 4 5 '  -0.01%  '

Set-TextCell 5 4 '325.31'
Set-TextCell 5 5 '  -0.09%  '

Set-TextCell 6 5 '  +0.03%  '

Set-TextCell 7 4 '0.4594'
Set-TextCell 7 5 '  -0.90%  '

Set-TextCell 8 4 '0.3825'
Set-TextCell 8 5 '  -1.21%  '

Set-TextCell 9 4 '0.07731'
Set-TextCell 9 5 '  -1.36%  '

Set-TextCell 10 5 '  +0.48%  '

Set-TextCell 11 4 '22.08'
Set-TextCell 11 5 '  -2.64%  '

Set-TextCell 12 4 '1.903.55'
Set-TextCell 12 5 '  -0.62%  '

Set-TextCell 13 5 '  -1.95%  '

Set-TextCell 14 4 '5.670'
Set-TextCell 14 5 '  -1.55%  '

Set-TextCell 15 4 '0.07027'
Set-TextCell 15 5 '  -0.27%  '

Set-TextCell 16 5 '  -0.06%  '

Set-TextCell 17 4 '83.91'
Set-TextCell 17 5 '  -3.30%  '

Set-TextCell 18 4 '0.000009470'
Set-TextCell 18 5 '  -3.64%  '

Set-TextCell 19 5 '  -2.28%  '

Set-TextCell 20 5 '  +0.00%  '

Set-TextCell 21 4 '28.919.25'
Set-TextCell 21 5 '  -1.69%  '

Set-TextCell 22 4 '5.321'
Set-TextCell 22 5 '  -2.71%  '

Set-TextCell 23 4 '10.88'
Set-TextCell 23 5 '  -1.87%  '

Set-TextCell 24 4 '2.093'
Set-TextCell 24 5 '  -0.27%  '

Set-TextCell 25 4 '158.50'
Set-TextCell 25 5 '  +0.74%  '

Set-TextCell 26 4 '19.06'
Set-TextCell 26 5 '  -1.59%  '

Set-TextCell 27 4 '5.674'
Set-TextCell 27 5 '  -1.53%  '

Set-TextCell 28 4 '117.53'
Set-TextCell 28 5 '  -0.98%  '

Set-TextCell 29 4 '1.857'
Set-TextCell 29 5 '  -0.30%  '

Set-TextCell 30 4 '0.09295'
Set-TextCell 30 5 '  -0.83%  '

Set-TextCell 31 4 '0.8692'
Set-TextCell 31 5 '  +1.06%  '

Set-TextCell 32 4 '5.082'
Set-TextCell 32 5 '  -2.09%  '

Set-TextCell 33 5 '  -4.11%  '

Set-TextCell 34 4 '3.140'
Set-TextCell 34 5 '  +0.48%  '

Set-TextCell 35 2 'Hedera'
Set-TextCell 35 3 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 35 4 '0.05725'
Set-TextCell 35 5 '  -0.81%  '

Set-TextCell 36 2 'TrustWalletToken'
Set-TextCell 36 3 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 36 4 '1.166'
Set-TextCell 36 5 '  +1.02%  '

Set-TextCell 37 4 '1.003'
Set-TextCell 37 5 '  +0.10%  '

Set-TextCell 38 4 '0.02044'
Set-TextCell 38 5 '  -1.89%  '

Set-TextCell 39 4 '0.5493'
Set-TextCell 39 5 '  -3.12%  '

Set-TextCell 40 4 '7.404'
Set-TextCell 40 5 '  -4.03%  '

Set-TextCell 41 4 '0.1755'
Set-TextCell 41 5 '  -1.58%  '

Set-TextCell 42 4 '2.856'
Set-TextCell 42 5 '  +4.70%  '

Set-TextCell 43 4 '9.310'
Set-TextCell 43 5 '  -1.21%  '

Set-TextCell 44 4 '0.5181'
Set-TextCell 44 5 '  -2.10%  '

Set-TextCell 45 2 'EnergySwap'
Set-TextCell 45 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 45 4 '11.21'
Set-TextCell 45 5 '  -2.83%  '

Set-TextCell 46 2 'Cronos'
Set-TextCell 46 3 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 46 4 '0.06914'
Set-TextCell 46 5 '  +0.58%  '

Set-TextCell 47 4 '2.099'
Set-TextCell 47 5 '  +0.75%  '

Set-TextCell 48 4 '0.000002575'
Set-TextCell 48 5 '  -10.09%  '

Set-TextCell 49 4 '1.779'
Set-TextCell 49 5 '  -2.11%  '

Set-TextCell 50 4 '110.54'
Set-TextCell 50 5 '  -0.71%  '

Set-TextCell 51 4 '0.2878'
Set-TextCell 51 5 '  -4.06%  '
